# Add three new per-market test-data worksheets (Netherlands, Austria,
# Denmark) at the end of the workbook, matching the structure/styling
# already used by the other market sheets (e.g. "Italy"), and make
# "Austria" the active/selected sheet — per commit message "Test data
# for Austria market added".

$wb = $excel.ActiveWorkbook

function New-CountrySheet($name) {
    # "Italy" has the exact column widths / row layout (no custom row
    # heights) that the newly added sheets use, so copy it as the
    # template and drop the copy at the very end of the tab strip.
    $template = $wb.Worksheets.Item("Italy")
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $lastSheet)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $name
    return $newSheet
}

$netherlands = New-CountrySheet "Netherlands"
$netherlands.Range("B4").Value = "NGC-3144/T2188/T2189/T2191"
$netherlands.Range("B2").Value = "Netherlands Market"
$netherlands.Range("C14").Select() | Out-Null

$austria = New-CountrySheet "Austria"
$austria.Range("B4").Value = "NGC-3817/T2295"
$austria.Range("B2").Value = "Austria Market"
$austria.Range("D15").Select() | Out-Null

$denmark = New-CountrySheet "Denmark"
$denmark.Range("B4").Value = "NGC-2913/T2783"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Select() | Out-Null

# Austria's "User Story" cell (A9) calls out the Fire Brigade Panel
# rather than the generic PR1D2 constant the template carries over.
$austria.Range("A9").Value = "Fire Brigade Panel"

# Leave Austria as the active/selected tab.
$austria.Activate()
